$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts existing Sentence Text/Length/Longest Word(s) columns to B/C/D)
$ws.Columns("A").Insert()
$ws.Columns("A").ColumnWidth = 40.16666666666667

# Header row
$ws.Range("A1").Value = "Test Description"

# Replace "Goose" with "MotherGoose" in the sample sentence / expected word columns
$ws.Range("B2").Value = "MotherGoose"
$ws.Range("C2").Value = 11
$ws.Range("D2").Value = "MotherGoose"

# New test-description column values (order matches shared-string insertion order)
$ws.Range("A2").Value = "one-word sentence"
$ws.Range("A5").Value = "longest word is followed by punctuation"
$ws.Range("A3").Value = "duplicate longest words"
$ws.Range("A4").Value = "single longest word"
$ws.Range("A6").Value = "longest word is apostrophized"
$ws.Range("A7").Value = "longest word is hyphenated"

$ws.Range("B10").Select() | Out-Null
